$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE: "001" -> "004" (must remain zero-padded text, not numeric 4)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"
$ws.Range("J2").ClearFormats()

# REPORT_DATE
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# TOTAL_ASSETS / FIXED_ASSET / MONETARYFUNDS
$ws.Range("O2").Value = 731743024.7
$ws.Range("P2").Value = 78747128.08
$ws.Range("Q2").Value = 15590515.52

# MONETARYFUNDS_RATIO no longer reported
$ws.Range("R2").ClearContents()

# ACCOUNTS_RECE
$ws.Range("S2").Value = 313250881.83

# ACCOUNTS_RECE_RATIO no longer reported
$ws.Range("T2").ClearContents()

# INVENTORY
$ws.Range("U2").Value = 124035973.77

# INVENTORY_RATIO no longer reported
$ws.Range("V2").ClearContents()

# TOTAL_LIABILITIES / ACCOUNTS_PAYABLE
$ws.Range("W2").Value = 167392722.68
$ws.Range("X2").Value = 94993856.88

# ADVANCE_RECEIVABLES / ADVANCE_RECEIVABLES_RATIO no longer reported
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()

# TOTAL_EQUITY
$ws.Range("AB2").Value = 564350302.02

# TOTAL_EQUITY_RATIO / TOTAL_ASSETS_RATIO / TOTAL_LIAB_RATIO no longer reported
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()

# CURRENT_RATIO / DEBT_ASSET_RATIO
$ws.Range("AF2").Value = 354.1355814212
$ws.Range("AG2").Value = 22.8758890799
